$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f1778345ee0>),
                (''model'',
                 AdaBoostClassifier(estimator=LGBMClassifier(boosting_type=''dart'',
                                                             class_weight=''balanced'',
                                                             colsample_bytree=0.9,
                                                             learning_rate=0.05,
                                                             max_depth=7,
                                                             num_leaves=2,
                                                             random_state=42,
                                                             subsample=0.9),
                                    random_state=42))])'
$ws.Range("B2").Value = 0.6111655011655011
$ws.Range("C2").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f1778249df0>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 50, ''model__estimator__subsample'': 0.9, ''model__estimator__num_leaves'': 2, ''model__estimator__min_child_samples'': 20, ''model__estimator__max_depth'': 7, ''model__estimator__learning_rate'': 0.05, ''model__estimator__colsample_bytree'': 0.9, ''model__estimator__class_weight'': ''balanced'', ''model__estimator__boosting_type'': ''dart''}'
$ws.Range("D2").Value = 0.962089050157198
$ws.Range("E2").Value = 0.5196241869241869
$ws.Range("F2").Value = 0.742857142857143
$ws.Range("G2").Value = 0.9712646964955758
$ws.Range("H2").Value = 0.5397650793650793
$ws.Range("I2").Value = 0.6842105263157895
$ws.Range("J2").Value = 0.9539148936170213
$ws.Range("K2").Value = 0.5226666666666666
$ws.Range("L2").Value = 0.8125
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("N2").Value = '[1 1 1 1 1 1 0 0 1 1 1 1 1 1 0 0 1 1 0 1 1 1 1 1]'
$ws.Range("O2").Value = 42

# Row 3
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f17782ae970>),
                (''model'',
                 AdaBoostClassifier(estimator=LGBMClassifier(boosting_type=''dart'',
                                                             class_weight=''balanced'',
                                                             colsample_bytree=0.7,
                                                             learning_rate=0.2,
                                                             max_depth=5,
                                                             min_child_samples=10,
                                                             num_leaves=10,
                                                             random_state=42,
                                                             subsample=0.7),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B3").Value = 0.6276301476301476
$ws.Range("C3").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f17782325e0>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__subsample'': 0.7, ''model__estimator__num_leaves'': 10, ''model__estimator__min_child_samples'': 10, ''model__estimator__max_depth'': 5, ''model__estimator__learning_rate'': 0.2, ''model__estimator__colsample_bytree'': 0.7, ''model__estimator__class_weight'': ''balanced'', ''model__estimator__boosting_type'': ''dart''}'
$ws.Range("D3").Value = 0.9540224852053723
$ws.Range("E3").Value = 0.4801240093240093
$ws.Range("F3").Value = 0.7647058823529411
$ws.Range("G3").Value = 0.9674092428847044
$ws.Range("H3").Value = 0.5702166666666666
$ws.Range("I3").Value = 0.7222222222222222
$ws.Range("J3").Value = 0.942340425531915
$ws.Range("K3").Value = 0.4383333333333333
$ws.Range("L3").Value = 0.8125
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range("N3").Value = '[0 1 1 1 1 0 1 1 1 1 1 0 0 0 1 1 1 1 1 1 1 0 1 1]'
$ws.Range("O3").Value = 69

# Row 4
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f17782ae5e0>),
                (''model'',
                 AdaBoostClassifier(estimator=LGBMClassifier(class_weight=''balanced'',
                                                             colsample_bytree=0.9,
                                                             learning_rate=0.01,
                                                             max_depth=3,
                                                             min_child_samples=5,
                                                             num_leaves=10,
                                                             random_state=42,
                                                             subsample=0.9),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B4").Value = 0.6393434343434342
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f17782408e0>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__subsample'': 0.9, ''model__estimator__num_leaves'': 10, ''model__estimator__min_child_samples'': 5, ''model__estimator__max_depth'': 3, ''model__estimator__learning_rate'': 0.01, ''model__estimator__colsample_bytree'': 0.9, ''model__estimator__class_weight'': ''balanced'', ''model__estimator__boosting_type'': ''gbdt''}'
$ws.Range("D4").Value = 0.9615375801503676
$ws.Range("E4").Value = 0.5035133755133755
$ws.Range("F4").Value = 0.5625
$ws.Range("G4").Value = 0.9669373468911572
$ws.Range("H4").Value = 0.5512738095238094
$ws.Range("I4").Value = 0.6923076923076923
$ws.Range("J4").Value = 0.9565777777777779
$ws.Range("K4").Value = 0.4904
$ws.Range("L4").Value = 0.4736842105263158
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N4").Value = '[0 1 1 1 0 1 0 1 1 0 0 1 0 1 0 1 0 1 0 0 1 1 0 1]'
$ws.Range("O4").Value = 23

# Row 5
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 AdaBoostClassifier(estimator=LGBMClassifier(class_weight=''balanced'',
                                                             colsample_bytree=0.5,
                                                             learning_rate=0.2,
                                                             max_depth=7,
                                                             num_leaves=20,
                                                             random_state=42,
                                                             subsample=0.9),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B5").Value = 0.6782051282051282
$ws.Range("C5").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__subsample'': 0.9, ''model__estimator__num_leaves'': 20, ''model__estimator__min_child_samples'': 20, ''model__estimator__max_depth'': 7, ''model__estimator__learning_rate'': 0.2, ''model__estimator__colsample_bytree'': 0.5, ''model__estimator__class_weight'': ''balanced'', ''model__estimator__boosting_type'': ''gbdt''}'
$ws.Range("D5").Value = 0.9758706873362157
$ws.Range("E5").Value = 0.558197557997558
$ws.Range("F5").Value = 0.5384615384615384
$ws.Range("G5").Value = 0.9845398871083151
$ws.Range("H5").Value = 0.5764515873015873
$ws.Range("I5").Value = 0.5833333333333334
$ws.Range("J5").Value = 0.9683265306122449
$ws.Range("K5").Value = 0.5630000000000001
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N5").Value = '[0 1 1 0 1 0 0 1 0 1 1 1 0 0 1 0 0 1 1 0 0 1 0 1]'
$ws.Range("O5").Value = 23

# Row 6
$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f1778232340>),
                (''model'',
                 AdaBoostClassifier(estimator=LGBMClassifier(class_weight=''balanced'',
                                                             colsample_bytree=0.5,
                                                             max_depth=5,
                                                             min_child_samples=10,
                                                             num_leaves=20,
                                                             random_state=42,
                                                             subsample=0.5),
                                    n_estimators=100, random_state=42))])'
$ws.Range("B6").Value = 0.7053446553446554
$ws.Range("C6").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f1778240d90>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 100, ''model__estimator__subsample'': 0.5, ''model__estimator__num_leaves'': 20, ''model__estimator__min_child_samples'': 10, ''model__estimator__max_depth'': 5, ''model__estimator__learning_rate'': 0.1, ''model__estimator__colsample_bytree'': 0.5, ''model__estimator__class_weight'': ''balanced'', ''model__estimator__boosting_type'': ''gbdt''}'
$ws.Range("D6").Value = 0.9646773873582138
$ws.Range("E6").Value = 0.5987598290598291
$ws.Range("F6").Value = 0.5600000000000001
$ws.Range("G6").Value = 0.9815272586860765
$ws.Range("H6").Value = 0.6106706349206349
$ws.Range("I6").Value = 0.5
$ws.Range("J6").Value = 0.9494999999999999
$ws.Range("K6").Value = 0.606
$ws.Range("L6").Value = 0.6363636363636364
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range("N6").Value = '[1 0 1 1 1 1 0 0 0 0 0 1 1 1 0 0 0 1 1 0 1 1 1 1]'
$ws.Range("O6").Value = 89
